$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.463.42"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "2.054.24"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  -0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "242.38"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.42%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.663"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("E7").Value = "  -0.01%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "54.24"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -5.51%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "58.24"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.68%  "
$ws.Range("E10").Value = "  -4.64%  "
$ws.Range("E11").Value = "  -2.10%  "
$ws.Range("E12").Value = "  -2.90%  "
$ws.Range("E13").Value = "  +2.63%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "14.69"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -5.25%  "
$ws.Range("D15").Value = "2.358.09"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("E16").Value = "  -4.89%  "
$ws.Range("D17").Value = "2.065.94"
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("D18").Value = "36.427.88"
$ws.Range("E18").Value = "  -1.28%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "16.78"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -7.31%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "71.86"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").Value = "0.0₃0855"
$ws.Range("E21").Value = "  -3.96%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "238.44"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("E23").Value = "  -2.80%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.04%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.35"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -3.55%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.31"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("E27").Value = "  -0.50%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "164.03"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.46%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "20.01"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.95%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.122"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.84%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.20"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +6.78%  "
$ws.Range("E32").Value = "  -7.55%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.45"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -5.85%  "
$ws.Range("E34").Value = "  -3.42%  "
$ws.Range("E35").Value = "  +0.01%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.82"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("E37").Value = "  -2.16%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0819"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -5.27%  "
$ws.Range("E39").Value = "  -5.24%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "4.85"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.28%  "
$ws.Range("E41").Value = "  -3.15%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.81"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -9.54%  "
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("E44").Value = "  -3.00%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0918"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -6.39%  "
$ws.Range("D46").Value = "1.397.65"
$ws.Range("E46").Value = "  +8.60%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "15.85"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -6.71%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.52"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +11.81%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("D51").Value = "2.246.29"
$ws.Range("E51").Value = "  +0.87%  "
